$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 862-863; everything currently at row 862
# downward shifts down by two rows (row 862 -> 864, ..., row 934 -> 936).
$ws.Rows("862:863").Insert()

# New weekly price entries (Fecha 2023-12-05 / serial 45265) for
# "Vega Monumental Concepción" - Limón, inserted at the top of this block.

# Row 862: 1a amarillo
$ws.Range("A862").Value = 11
$ws.Range("B862").Value = "Vega Monumental Concepción"
$ws.Range("C862").Value = "Bíobío"
$ws.Range("D862").Value = 45265
$ws.Range("E862").Value = 8
$ws.Range("F862").Value = "Fruta"
$ws.Range("G862").Value = 100102
$ws.Range("H862").Value = "Cítricos"
$ws.Range("I862").Value = 100102003
$ws.Range("J862").Value = "Limón"
$ws.Range("K862").Value = "Sin especificar"
$ws.Range("L862").Value = "1a amarillo"
$ws.Range("M862").Value = 150
$ws.Range("N862").Value = 12000
$ws.Range("O862").Value = 12000
$ws.Range("P862").Value = 12000
$ws.Range("Q862").Value = "$/malla 16 kilos"
$ws.Range("R862").Value = "Región de O'Higgins"
$ws.Range("S862").Value = 750
$ws.Range("T862").Value = 16

# Row 863: 1a plateado
$ws.Range("A863").Value = 11
$ws.Range("B863").Value = "Vega Monumental Concepción"
$ws.Range("C863").Value = "Bíobío"
$ws.Range("D863").Value = 45265
$ws.Range("E863").Value = 8
$ws.Range("F863").Value = "Fruta"
$ws.Range("G863").Value = 100102
$ws.Range("H863").Value = "Cítricos"
$ws.Range("I863").Value = 100102003
$ws.Range("J863").Value = "Limón"
$ws.Range("K863").Value = "Sin especificar"
$ws.Range("L863").Value = "1a plateado"
$ws.Range("M863").Value = 120
$ws.Range("N863").Value = 14000
$ws.Range("O863").Value = 14000
$ws.Range("P863").Value = 14000
$ws.Range("Q863").Value = "$/malla 16 kilos"
$ws.Range("R863").Value = "Región de O'Higgins"
$ws.Range("S863").Value = 875
$ws.Range("T863").Value = 16
